# #12 Colors added on web page on solution slides
# Target slide: "Web page" (solution) slide -> slide 16 in this deck.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

# ---------------------------------------------------------------------
# 1) Rectangle 3 (shape 1): give the rectangle a dark (tx1) fill and an
#    outline of the same colour.
# ---------------------------------------------------------------------
$rect = $s.Shapes.Item(1)
$rect.Fill.Visible = 1
$rect.Fill.ForeColor.SchemeColor = "tx1"
$rect.Line.Visible = 1
$rect.Line.ForeColor.RGB = 0

# ---------------------------------------------------------------------
# 2) Body placeholder (shape 3): bold + colour a couple of words.
# ---------------------------------------------------------------------
$body = $s.Shapes.Item(3)
$tr = $body.TextFrame.TextRange

# "web page " (chars 5-13 of "The web page will ...") -> bold + dark colour
$webPage = $tr.Characters(5, 9)
$webPage.Font.Bold = 1
$webPage.Font.Color.RGB = 0

# "strategies" -> bold + blue (0070C0)
$strategies = $tr.Characters(35, 10)
$strategies.Font.Bold = 1
$strategies.Font.Color.RGB = 0x00 + 0x70 * 256 + 0xC0 * 65536

# ---------------------------------------------------------------------
# 3) Picture "Graphique 4" (shape 4): nudge it up a little.
# ---------------------------------------------------------------------
$pic = $s.Shapes.Item(4)
$pic.Top = 2989719 / 12700
